# Fruta / hortaliza, semanal
#
# The weekly update inserts two new price records (dated 2022-01-24, serial
# 44585) into the "Macroferia Regional de Talca - Plátano" sheet. The new
# rows are inserted right after the existing row 477 (pushing the former
# rows 478-508 down to 480-510), and the sheet dimension grows from
# A1:T508 to A1:T510.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 478-479; everything that was in 478:508 shifts
# down to 480:510 automatically (formats follow the row below, matching
# the D-column date style already used throughout the sheet).
$ws.Rows("478:479").Insert()

# New row 478: Pintón, 450 cajas, 12000/12000/12000, 600 $/Kg
$ws.Range("A478").Value = 5
$ws.Range("B478").Value = "Macroferia Regional de Talca"
$ws.Range("C478").Value = "Maule"
$ws.Range("D478").Value = 44585
$ws.Range("D478").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E478").Value = 7
$ws.Range("F478").Value = "Fruta"
$ws.Range("G478").Value = 100108
$ws.Range("H478").Value = "Tropicales y subtropicales"
$ws.Range("I478").Value = 100108006
$ws.Range("J478").Value = "Plátano"
$ws.Range("K478").Value = "Sin especificar"
$ws.Range("L478").Value = "Pintón"
$ws.Range("M478").Value = 450
$ws.Range("N478").Value = 12000
$ws.Range("O478").Value = 12000
$ws.Range("P478").Value = 12000
$ws.Range("Q478").Value = "`$/caja 20 kilos"
$ws.Range("R478").Value = "Ecuador"
$ws.Range("S478").Value = 600
$ws.Range("T478").Value = 20

# New row 479: Primera Pintón, 300 cajas, 13000/13000/13000, 650 $/Kg
$ws.Range("A479").Value = 5
$ws.Range("B479").Value = "Macroferia Regional de Talca"
$ws.Range("C479").Value = "Maule"
$ws.Range("D479").Value = 44585
$ws.Range("D479").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E479").Value = 7
$ws.Range("F479").Value = "Fruta"
$ws.Range("G479").Value = 100108
$ws.Range("H479").Value = "Tropicales y subtropicales"
$ws.Range("I479").Value = 100108006
$ws.Range("J479").Value = "Plátano"
$ws.Range("K479").Value = "Sin especificar"
$ws.Range("L479").Value = "Primera Pintón"
$ws.Range("M479").Value = 300
$ws.Range("N479").Value = 13000
$ws.Range("O479").Value = 13000
$ws.Range("P479").Value = 13000
$ws.Range("Q479").Value = "`$/caja 20 kilos"
$ws.Range("R479").Value = "Ecuador"
$ws.Range("S479").Value = 650
$ws.Range("T479").Value = 20
